$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new "01-sep" column (BK) to the daily tracker, one row per entry.
$ws.Range("BK1").Value = "01-sep"
$ws.Range("BK2").Value = 0
$ws.Range("BK3").Value = 27.239267029014986
$ws.Range("BK4").Value = 23.297642305215735
$ws.Range("BK5").Value = 17.640707357263107
$ws.Range("BK6").Value = 0
$ws.Range("BK7").Value = 19.003421510140669
$ws.Range("BK8").Value = 11.551901808677695
$ws.Range("BK9").Value = 9.6931865424973367
$ws.Range("BK10").Value = 27.977605772145317
$ws.Range("BK11").Value = 18.010741131289311
$ws.Range("BK12").Value = 0
$ws.Range("BK13").Value = 11.226310241047724
$ws.Range("BK14").Value = 0
$ws.Range("BK15").Value = 0
$ws.Range("BK16").Value = 24.713057225738531
$ws.Range("BK17").Value = 0
$ws.Range("BK18").Value = 0

# Match the saved selection state left behind by the author's session.
[void]$ws.Range("BM4").Select()
